# Swap the data between row 25 and row 26 on the "Artfynd" sheet.
# (Two species observations - "Knärot" and "Kandelabersvamp" - had their
#  row positions exchanged; most columns keep their values but columns
#  A,B,D,E,F,G,H,I,J,K,L,Q,R,AC move from one row to the other.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 1000

# Columns whose *content* needs to be exchanged between row 25 and row 26.
$cols = @("A","B","D","E","F","G","H","I","J","Q","R","AC")

foreach ($col in $cols) {
    $src = $col + "25"
    $dst = $col + "26"
    $tmp = $col + $scratchRow

    # row25 -> scratch
    $ws.Range($src).Copy()
    $ws.Range($tmp).PasteSpecial()

    # row26 -> row25
    $ws.Range($dst).Copy()
    $ws.Range($src).PasteSpecial()

    # scratch -> row26
    $ws.Range($tmp).Copy()
    $ws.Range($dst).PasteSpecial()

    # tidy up the scratch cell
    $ws.Range($tmp).Clear()
}

# Column K: "Åldersstadium" text ("blomning" on row 25, blank on row 26)
# also exchanges between the two rows.
$ws.Range("K26").Copy()
$ws.Range("K" + $scratchRow).PasteSpecial()
$ws.Range("K26").Value = "'blomning"
$ws.Range("K25").Clear()

# Column L ("Kön") only existed on row 25 (empty placeholder); after the
# swap it belongs conceptually to row 26 instead.
$ws.Range("L25").Clear()

$ws.Range("K" + $scratchRow).Clear()
